$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para1 = "The resulting output should be a set of symbols mirroring the originals – to a degree."
$para2 = "Due to time constraints, I was limited on time to train the model. It’s taken a day so far to get this running through all of the epochs."
$para3 = "That being said, due to the lowered batch size and epoch count I had to give my program due to that time constraint, I estimate that the symbols may not resemble their original counterparts as well as I’d hope."
$para4a = "Currently, the program is running smoothly, and I await the results of its generation "
$para4b = "to report them later."

$tr.Text = $para1 + [char]13 + $para2 + [char]13 + $para3 + [char]13 + $para4a + $para4b

$p4 = $tr.Paragraphs(4)
$splitStart = $p4.Start + $para4a.Length
$splitLen = $para4b.Length
$sub = $tr.Characters($splitStart, $splitLen)
$sub.Text = $para4b
